# One-click update from Van Paper 07:23 AM on 2025-11-13
#
# Business changes to the Sheet1 leaderboard table:
#   - Row 28 "NICOLLET COURT RETAIL MALL" is renamed to "NICOLLET RETAIL LLC"
#     and its Salesperson code moves from 015 to 023.
#   - A brand new customer "MAYNARD'S" (Salesperson 040, customer # 0008369)
#     is inserted immediately after it, pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing customer and update its salesperson code.
$ws.Range("A28").Value = "NICOLLET RETAIL LLC"
$ws.Range("C28").Value = "023"

# Make room for the new customer row by inserting a row at 29 - this shifts
# the old rows 29 (HOLY FAMILY MARONITE CHURCH) and 30 (SCHMITT MUSIC CTR)
# down to 30 and 31 respectively.
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).RowHeight = 13.05

# Populate the newly inserted row with the new customer's data.
$ws.Range("A29").Value = "MAYNARD'S"
$ws.Range("B29").Value = "Norman, Ryan M"
$ws.Range("C29").Value = "040"
$ws.Range("E29").Value = "0008369"
